$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A5").Value = "2/1/2/108"
$ws.Range("C5").Value = "1 tỉ"

$ws.Range("D7").Select()
